$wb = $excel.ActiveWorkbook

# Sheet ALC (line 2476 in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 186.66667
$ws.Range("I38").Value = 186.66667
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 560.00001
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -188.00001

# Sheet ALC (line 6159 in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2083.5
$ws.Range("J112").Value = 2083.5
$ws.Range("L112").Value = 6250.5
$ws.Range("N112").Value = -8466.5

# Sheet ALC (line 7139 in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6435.303
$ws.Range("I132").Value = 7303.44
$ws.Range("K132").Value = 21910.32
$ws.Range("M132").Value = -19380.32

# Sheet ALC (line 7390 in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1167598.8
$ws.Range("I137").Value = 2778997.2
$ws.Range("K137").Value = 8336991.600000001
$ws.Range("M137").Value = -8334441.600000001

# Sheet ALC (line 7442 in diff)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4044.8215
$ws.Range("I138").Value = 5539.3335
$ws.Range("K138").Value = 16618.0005
$ws.Range("M138").Value = -11478.0005

# Sheet ARM (line 9788 in diff)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 23024.85
$ws.Range("I45").Value = 31132.572
$ws.Range("J45").Value = 4106.8335
$ws.Range("K45").Value = 31132.572
$ws.Range("L45").Value = 4106.8335
$ws.Range("M45").Value = -30755.572
$ws.Range("N45").Value = -4860.8335

# Sheet ARM (line 10612 in diff)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0

# Sheet ARM (line 10756 in diff)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

# Sheet ARM (line 13985 in diff)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3067.7576
$ws.Range("I132").Value = 2586
$ws.Range("K132").Value = 7758
$ws.Range("M132").Value = -5228

# Sheet BSM (line 15442 in diff)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20008116
$ws.Range("I20").Value = 29421572
$ws.Range("K20").Value = 29421572
$ws.Range("M20").Value = -29421325

# Sheet BSM (line 19287 in diff)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 46273.168
$ws.Range("J100").Value = 46273.168
$ws.Range("L100").Value = 46273.168
$ws.Range("N100").Value = -48437.168

# Sheet BSM (line 19529 in diff)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13686485
$ws.Range("I105").Value = 770970.6
$ws.Range("J105").Value = 41670100
$ws.Range("K105").Value = 770970.6
$ws.Range("L105").Value = 41670100
$ws.Range("M105").Value = -769223.6
$ws.Range("N105").Value = -41673594

# Sheet BSM (line 20899 in diff)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2905.3818
$ws.Range("I134").Value = 2568.762
$ws.Range("J134").Value = 3992.923
$ws.Range("K134").Value = 7706.286
$ws.Range("L134").Value = 11978.769
$ws.Range("M134").Value = -5171.286
$ws.Range("N134").Value = -17048.769

# Sheet CRP (line 21624 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1872.5
$ws.Range("I7").Value = 1750
$ws.Range("K7").Value = 1750
$ws.Range("M7").Value = -1637

# Sheet CRP (line 22794 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3319.1956
$ws.Range("I31").Value = 2364.4
$ws.Range("K31").Value = 2364.4
$ws.Range("M31").Value = -2069.4

# Sheet CRP (line 22941 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3319.1956
$ws.Range("I34").Value = 2364.4
$ws.Range("K34").Value = 2364.4
$ws.Range("M34").Value = -2162.4

# Sheet CRP (line 25435 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 22398.666
$ws.Range("I86").Value = 28631.334
$ws.Range("J86").Value = 9933.333000000001
$ws.Range("K86").Value = 28631.334
$ws.Range("L86").Value = 9933.333000000001
$ws.Range("M86").Value = -27508.334
$ws.Range("N86").Value = -12179.333

# Sheet CRP (line 25582 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 22398.666
$ws.Range("I89").Value = 28631.334
$ws.Range("J89").Value = 9933.333000000001
$ws.Range("K89").Value = 143156.67
$ws.Range("L89").Value = 49666.665
$ws.Range("M89").Value = -137540.67
$ws.Range("N89").Value = -60898.665

# Sheet CRP (line 25931 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 6271.25
$ws.Range("J96").Value = 6271.25
$ws.Range("L96").Value = 6271.25
$ws.Range("N96").Value = -11763.25

# Sheet CRP (line 25980 in diff)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 43084.5
$ws.Range("J97").Value = 49990
$ws.Range("L97").Value = 49990
$ws.Range("N97").Value = -51972

# Sheet CUL (line 33088 in diff)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 993.2
$ws.Range("I98").Value = 513.2
$ws.Range("J98").Value = 1473.2
$ws.Range("K98").Value = 1539.6
$ws.Range("L98").Value = 4419.6
$ws.Range("M98").Value = -41.60000000000014
$ws.Range("N98").Value = -7415.6

# Sheet GSM (line 35879 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("N12").Value = 0

# Sheet GSM (line 37319 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 16285.571
$ws.Range("J42").Value = 16285.571
$ws.Range("L42").Value = 16285.571
$ws.Range("N42").Value = -17255.571

# Sheet GSM (line 39130 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111114970
$ws.Range("I80").Value = 250002990
$ws.Range("J80").Value = 4549.6
$ws.Range("K80").Value = 250002990
$ws.Range("L80").Value = 4549.6
$ws.Range("M80").Value = -250001992
$ws.Range("N80").Value = -6545.6

# Sheet GSM (line 39280 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 111114970
$ws.Range("I83").Value = 250002990
$ws.Range("J83").Value = 4549.6
$ws.Range("K83").Value = 1250014950
$ws.Range("L83").Value = 22748
$ws.Range("M83").Value = -1250009958
$ws.Range("N83").Value = -32732

# Sheet GSM (line 40818 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 16285.571
$ws.Range("J115").Value = 16285.571
$ws.Range("L115").Value = 16285.571
$ws.Range("N115").Value = -18635.571

# Sheet GSM (line 41624 in diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3651.0588
$ws.Range("I132").Value = 2952.2222
$ws.Range("K132").Value = 8856.6666
$ws.Range("M132").Value = -6326.6666

# Sheet LTW (line 46565 in diff)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# Sheet LTW (line 46902 in diff)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1499.5
$ws.Range("I100").Value = 1366.3334
$ws.Range("J100").Value = 1632.6666
$ws.Range("K100").Value = 1366.3334
$ws.Range("L100").Value = 1632.6666
$ws.Range("M100").Value = -825.3334
$ws.Range("N100").Value = -2714.6666

# Sheet WVR (line 49070 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Sheet WVR (line 49505 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 5000
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4858

# Sheet WVR (line 51117 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 29999
$ws.Range("J46").Value = 29999
$ws.Range("L46").Value = 29999
$ws.Range("N46").Value = -30461

# Sheet WVR (line 53932 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 22368.4
$ws.Range("J104").Value = 22368.4
$ws.Range("L104").Value = 22368.4
$ws.Range("N104").Value = -29356.4

# Sheet WVR (line 55001 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3475.5
$ws.Range("I126").Value = 3493.5
$ws.Range("K126").Value = 10480.5
$ws.Range("M126").Value = -8010.5

# Sheet WVR (line 55384 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 29999
$ws.Range("J134").Value = 29999
$ws.Range("L134").Value = 89997
$ws.Range("N134").Value = -95067

# Sheet WVR (line 55479 in diff)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10313.625
$ws.Range("I136").Value = 10266.265
$ws.Range("J136").Value = 10428.643
$ws.Range("K136").Value = 30798.795
$ws.Range("L136").Value = 31285.929
$ws.Range("M136").Value = -28248.795
